$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateNewCompany")

# Update the city/state cell value (E2) from "New York City, NY" to "New York City"
$ws.Range("E2").Value = "New York City"

# Select the sheet and the edited cell, matching the resulting view state
$ws.Activate()
$ws.Range("E2").Select()
